# Refresh the crypto table with the latest scrape: columns D (Price) and
# E (Volume(1h) change) are plain text cells, so values are written as
# strings to preserve their exact display (trailing zeros, "NN.NNN.NN"
# thousands-style separators, padded "  +x.xx%  " percentages, etc.).
# A handful of new D-column values read as plain decimals (e.g. "233.61");
# those get a leading apostrophe so Excel keeps storing them as text
# instead of silently converting them to a Double.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "37.704.58"
$ws.Range("E2").Value = '  -0.17%  '

$ws.Range("D3").Value = "2.076.89"
$ws.Range("E3").Value = '  -1.76%  '

$ws.Range("E4").Value = '  +0.03%  '

$ws.Range("D5").Value = "'233.61"
$ws.Range("E5").Value = '  -0.81%  '

$ws.Range("D6").Value = "'0.624"
$ws.Range("E6").Value = '  -0.33%  '

$ws.Range("E7").Value = '  -0.04%  '

$ws.Range("D8").Value = "'58.22"
$ws.Range("E8").Value = '  -0.14%  '

$ws.Range("D9").Value = "'0.391"
$ws.Range("E9").Value = '  -0.25%  '

$ws.Range("E10").Value = '  +0.23%  '

$ws.Range("D11").Value = "'0.106"
$ws.Range("E11").Value = '  +2.73%  '

$ws.Range("D12").Value = "2.382.22"
$ws.Range("E12").Value = '  -1.96%  '

$ws.Range("D13").Value = "'14.83"
$ws.Range("E13").Value = '  +1.40%  '

$ws.Range("D14").Value = "'20.87"
$ws.Range("E14").Value = '  -1.28%  '

$ws.Range("D15").Value = "'0.771"
$ws.Range("E15").Value = '  -2.10%  '

$ws.Range("E16").Value = '  +1.09%  '

$ws.Range("D17").Value = "2.082.06"
$ws.Range("E17").Value = '  -1.57%  '

$ws.Range("D18").Value = "37.572.99"
$ws.Range("E18").Value = '  -0.26%  '

$ws.Range("D19").Value = "'6.19"
$ws.Range("E19").Value = '  -0.33%  '

$ws.Range("D20").Value = "'71.05"
$ws.Range("E20").Value = '  +1.27%  '

$ws.Range("D21").Value = "0.0₃0832"
$ws.Range("E21").Value = '  +1.04%  '

$ws.Range("D22").Value = "'227.67"
$ws.Range("E22").Value = '  +0.03%  '

$ws.Range("E23").Value = '  -0.06%  '

$ws.Range("E24").Value = '  -0.27%  '

$ws.Range("D25").Value = "'2.39"
$ws.Range("E25").Value = '  -1.30%  '

$ws.Range("D26").Value = "'169.41"
$ws.Range("E26").Value = '  +0.79%  '

$ws.Range("E27").Value = '  +3.15%  '

$ws.Range("E28").Value = '  -0.07%  '

$ws.Range("E29").Value = '  -0.31%  '

$ws.Range("E30").Value = '  -2.30%  '

$ws.Range("E31").Value = '  +2.22%  '

$ws.Range("D32").Value = "'4.67"
$ws.Range("E32").Value = '  +0.61%  '

$ws.Range("D33").Value = "'0.0629"
$ws.Range("E33").Value = '  +1.18%  '

$ws.Range("E34").Value = '  +1.16%  '

$ws.Range("E35").Value = '  -3.87%  '

$ws.Range("E36").Value = '  +2.75%  '

$ws.Range("E37").Value = '  -2.87%  '

$ws.Range("D38").Value = "'0.999"
$ws.Range("E38").Value = '  -0.03%  '

$ws.Range("D39").Value = "'5.34"
$ws.Range("E39").Value = '  -5.50%  '

$ws.Range("D40").Value = "'0.0978"
$ws.Range("E40").Value = '  +1.41%  '

$ws.Range("D41").Value = "'97.97"
$ws.Range("E41").Value = '  +0.19%  '

$ws.Range("E42").Value = '  +0.48%  '

$ws.Range("E43").Value = '  -2.77%  '

$ws.Range("D44").Value = "1.452.14"
$ws.Range("E44").Value = '  -1.66%  '

$ws.Range("E45").Value = '  -0.52%  '

$ws.Range("D46").Value = "'16.48"
$ws.Range("E46").Value = '  +5.50%  '

$ws.Range("D47").Value = "'4.23"
$ws.Range("E47").Value = '  +0.31%  '

$ws.Range("E48").Value = '  +0.76%  '

$ws.Range("E49").Value = '  +0.57%  '

$ws.Range("D50").Value = "'3.01"
$ws.Range("E50").Value = '  -0.92%  '

$ws.Range("D51").Value = "2.265.91"
$ws.Range("E51").Value = '  -1.84%  '
